$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.426.08'
$ws.Range("D2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.854.34'
$ws.Range("D3").ClearFormats()

$ws.Range("E3").Value = '  +0.36%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9994'
$ws.Range("D4").ClearFormats()

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.35'
$ws.Range("D5").ClearFormats()

$ws.Range("E5").Value = '  +0.23%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6324'
$ws.Range("D6").ClearFormats()

$ws.Range("E6").Value = '  +0.84%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07581'
$ws.Range("D8").ClearFormats()

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2925'
$ws.Range("D9").ClearFormats()

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.61'
$ws.Range("D10").ClearFormats()

$ws.Range("E10").Value = '  -0.85%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07753'
$ws.Range("D11").ClearFormats()

$ws.Range("E11").Value = '  +0.10%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.854.79'
$ws.Range("D12").ClearFormats()

$ws.Range("E12").Value = '  +0.38%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.036'
$ws.Range("D13").ClearFormats()

$ws.Range("E13").Value = '  +0.13%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6853'
$ws.Range("D14").ClearFormats()

$ws.Range("E14").Value = '  +0.70%  '

$ws.Range("E15").Value = '  -2.73%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.28'
$ws.Range("D16").ClearFormats()

$ws.Range("E16").Value = '  -0.14%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.115.28'
$ws.Range("D17").ClearFormats()

$ws.Range("E17").Value = '  +0.80%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.148'
$ws.Range("D18").ClearFormats()

$ws.Range("E18").Value = '  -0.40%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '29.435.41'
$ws.Range("D19").ClearFormats()

$ws.Range("E19").Value = '  -0.02%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '230.37'
$ws.Range("D20").ClearFormats()

$ws.Range("E20").Value = '  +1.01%  '

$ws.Range("E21").Value = '  +0.09%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("D22").ClearFormats()

$ws.Range("E22").Value = '  +0.00%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.527'
$ws.Range("D23").ClearFormats()

$ws.Range("E23").Value = '  +1.46%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9998'
$ws.Range("D24").ClearFormats()

$ws.Range("E24").Value = '  -0.08%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.27'
$ws.Range("D25").ClearFormats()

$ws.Range("E25").Value = '  +0.78%  '

$ws.Range("E26").Value = '  +1.91%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.470'
$ws.Range("D27").ClearFormats()

$ws.Range("E27").Value = '  +0.77%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.74'
$ws.Range("D28").ClearFormats()

$ws.Range("E28").Value = '  +0.28%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.414'
$ws.Range("D29").ClearFormats()

$ws.Range("E29").Value = '  +4.88%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.481'
$ws.Range("D30").ClearFormats()

$ws.Range("E30").Value = '  +1.43%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05695'
$ws.Range("D31").ClearFormats()

$ws.Range("E31").Value = '  +0.57%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.154'
$ws.Range("D32").ClearFormats()

$ws.Range("E32").Value = '  +0.82%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.061'
$ws.Range("D33").ClearFormats()

$ws.Range("E33").Value = '  +0.82%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.830'
$ws.Range("D34").ClearFormats()

$ws.Range("E34").Value = '  -0.66%  '

$ws.Range("E35").Value = '  -0.36%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7000'
$ws.Range("D36").ClearFormats()

$ws.Range("E36").Value = '  -0.21%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.588'
$ws.Range("D37").ClearFormats()

$ws.Range("E37").Value = '  +0.00%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.250.59'
$ws.Range("D38").ClearFormats()

$ws.Range("E38").Value = '  +2.10%  '

$ws.Range("E39").Value = '  +2.29%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.772'
$ws.Range("D40").ClearFormats()

$ws.Range("E40").Value = '  +0.39%  '

$ws.Range("E41").Value = '  -0.37%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9084'
$ws.Range("D42").ClearFormats()

$ws.Range("E42").Value = '  +0.33%  '

$ws.Range("E43").Value = '  -0.01%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.018.21'
$ws.Range("D44").ClearFormats()

$ws.Range("E44").Value = '  +0.58%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.59'
$ws.Range("D45").ClearFormats()

$ws.Range("E45").Value = '  -0.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '66.16'
$ws.Range("D46").ClearFormats()

$ws.Range("E46").Value = '  +0.22%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.153'
$ws.Range("D47").ClearFormats()

$ws.Range("E47").Value = '  -0.01%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000118'
$ws.Range("D48").ClearFormats()

$ws.Range("E48").Value = '  -2.77%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1166'
$ws.Range("D49").ClearFormats()

$ws.Range("E49").Value = '  +1.02%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.059'
$ws.Range("D50").ClearFormats()

$ws.Range("E50").Value = '  +0.41%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3971'
$ws.Range("D51").ClearFormats()

$ws.Range("E51").Value = '  -1.11%  '
